$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 373.55554
$ws.Range("I6").Value = 80.14286
$ws.Range("J6").Value = 1400.5
$ws.Range("K6").Value = 240.42858
$ws.Range("L6").Value = 4201.5
$ws.Range("M6").Value = -128.42858
$ws.Range("N6").Value = -4425.5
$ws.Range("H29").Value = 6767.6665
$ws.Range("I29").Value = 503
$ws.Range("K29").Value = 1509
$ws.Range("M29").Value = -1228
$ws.Range("H33").Value = 113.5
$ws.Range("I33").Value = 96.25
$ws.Range("J33").Value = 125
$ws.Range("K33").Value = 96.25
$ws.Range("L33").Value = 125
$ws.Range("M33").Value = 132.75
$ws.Range("N33").Value = -583
$ws.Range("H45").Value = 1933.3334
$ws.Range("I45").Value = 1900
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 5700
$ws.Range("L45").Value = 6000
$ws.Range("M45").Value = -5508
$ws.Range("N45").Value = -6384
$ws.Range("H62").Value = 2999.1667
$ws.Range("I62").Value = 2998.3333
$ws.Range("K62").Value = 2998.3333
$ws.Range("M62").Value = -2374.3333
$ws.Range("H65").Value = 2999.1667
$ws.Range("I65").Value = 2998.3333
$ws.Range("K65").Value = 14991.6665
$ws.Range("M65").Value = -11871.6665
$ws.Range("H86").Value = 1570.7142
$ws.Range("I86").Value = 1499.1666
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 1499.1666
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -376.1666
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 1570.7142
$ws.Range("I89").Value = 1499.1666
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 7495.833000000001
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -1879.833000000001
$ws.Range("N89").Value = -21232
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("H113").Value = 10947.417
$ws.Range("I113").Value = 12764.9
$ws.Range("K113").Value = 12764.9
$ws.Range("M113").Value = -9510.9
$ws.Range("H125").Value = 1493.5
$ws.Range("I125").Value = 1566.6666
$ws.Range("J125").Value = 1420.3334
$ws.Range("K125").Value = 14099.9994
$ws.Range("L125").Value = 12783.0006
$ws.Range("M125").Value = -11639.9994
$ws.Range("N125").Value = -17703.0006
$ws.Range("H129").Value = 861.4194
$ws.Range("J129").Value = 878.9107
$ws.Range("L129").Value = 2636.7321
$ws.Range("N129").Value = -12636.7321
$ws.Range("H132").Value = 938.85187
$ws.Range("I132").Value = 799.8570999999999
$ws.Range("J132").Value = 2301
$ws.Range("K132").Value = 2399.5713
$ws.Range("L132").Value = 6903
$ws.Range("M132").Value = 130.4287000000004
$ws.Range("N132").Value = -11963
$ws.Range("H135").Value = 514
$ws.Range("I135").Value = 498.94736
$ws.Range("K135").Value = 4490.52624
$ws.Range("M135").Value = -1955.52624
$ws.Range("H141").Value = 968132.4399999999
$ws.Range("I141").Value = 1219181.8
$ws.Range("K141").Value = 3657545.4
$ws.Range("M141").Value = -3652365.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2813.4683
$ws.Range("I32").Value = 2106.5693
$ws.Range("K32").Value = 2106.5693
$ws.Range("M32").Value = -1819.5693
$ws.Range("H61").Value = 1958.1333
$ws.Range("I61").Value = 1117.5652
$ws.Range("J61").Value = 4720
$ws.Range("K61").Value = 1117.5652
$ws.Range("L61").Value = 4720
$ws.Range("M61").Value = -905.5652
$ws.Range("N61").Value = -5144
$ws.Range("H110").Value = 1992.1
$ws.Range("I110").Value = 1130.1428
$ws.Range("J110").Value = 4003.3333
$ws.Range("K110").Value = 1130.1428
$ws.Range("L110").Value = 4003.3333
$ws.Range("M110").Value = 914.8571999999999
$ws.Range("N110").Value = -8093.3333
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H136").Value = 1958.1333
$ws.Range("I136").Value = 1117.5652
$ws.Range("J136").Value = 4720
$ws.Range("K136").Value = 3352.6956
$ws.Range("L136").Value = 14160
$ws.Range("M136").Value = -802.6956
$ws.Range("N136").Value = -19260

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1806.7894
$ws.Range("I20").Value = 1633.9333
$ws.Range("J20").Value = 2455
$ws.Range("K20").Value = 1633.9333
$ws.Range("L20").Value = 2455
$ws.Range("M20").Value = -1386.9333
$ws.Range("N20").Value = -2949
$ws.Range("H75").Value = 11867
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 11867
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H94").Value = 583.0357
$ws.Range("I94").Value = 420.2
$ws.Range("K94").Value = 420.2
$ws.Range("M94").Value = 30.80000000000001
$ws.Range("H99").Value = 1397.8572
$ws.Range("I99").Value = 1157
$ws.Range("K99").Value = 1157
$ws.Range("M99").Value = 341
$ws.Range("H105").Value = 2283.8262
$ws.Range("I105").Value = 2239.476
$ws.Range("K105").Value = 2239.476
$ws.Range("M105").Value = -492.4760000000001
$ws.Range("H108").Value = 94981
$ws.Range("J108").Value = 94981
$ws.Range("L108").Value = 94981
$ws.Range("N108").Value = -102661

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2326.7273
$ws.Range("I31").Value = 1899.8
$ws.Range("J31").Value = 2682.5
$ws.Range("K31").Value = 1899.8
$ws.Range("L31").Value = 2682.5
$ws.Range("M31").Value = -1604.8
$ws.Range("N31").Value = -3272.5
$ws.Range("H34").Value = 2326.7273
$ws.Range("I34").Value = 1899.8
$ws.Range("J34").Value = 2682.5
$ws.Range("K34").Value = 1899.8
$ws.Range("L34").Value = 2682.5
$ws.Range("M34").Value = -1697.8
$ws.Range("N34").Value = -3086.5
$ws.Range("H99").Value = 2155.3635
$ws.Range("I99").Value = 1451.6666
$ws.Range("J99").Value = 2999.8
$ws.Range("K99").Value = 1451.6666
$ws.Range("L99").Value = 2999.8
$ws.Range("M99").Value = 46.33339999999998
$ws.Range("N99").Value = -5995.8
$ws.Range("H107").Value = 419.15384
$ws.Range("J107").Value = 696
$ws.Range("L107").Value = 696
$ws.Range("N107").Value = -4536
$ws.Range("H122").Value = 4530.273
$ws.Range("I122").Value = 2856
$ws.Range("K122").Value = 8568
$ws.Range("M122").Value = -6118
$ws.Range("H125").Value = 29997.5
$ws.Range("J125").Value = 29997.5
$ws.Range("L125").Value = 29997.5
$ws.Range("N125").Value = -34917.5
$ws.Range("H126").Value = 2155.3635
$ws.Range("I126").Value = 1451.6666
$ws.Range("J126").Value = 2999.8
$ws.Range("K126").Value = 4354.9998
$ws.Range("L126").Value = 8999.400000000001
$ws.Range("M126").Value = -1884.9998
$ws.Range("N126").Value = -13939.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 141.5
$ws.Range("I14").Value = 141.5
$ws.Range("K14").Value = 424.5
$ws.Range("M14").Value = -251.5
$ws.Range("H92").Value = 300
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H131").Value = 814.74
$ws.Range("J131").Value = 814.74
$ws.Range("L131").Value = 2444.22
$ws.Range("N131").Value = -12524.22
$ws.Range("H132").Value = 1989.8
$ws.Range("I132").Value = 1499.6666
$ws.Range("K132").Value = 13496.9994
$ws.Range("M132").Value = -10966.9994
$ws.Range("H140").Value = 1492.5
$ws.Range("I140").Value = 835.45
$ws.Range("K140").Value = 2506.35
$ws.Range("M140").Value = 2673.65

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3225.4614
$ws.Range("I102").Value = 3470.111
$ws.Range("K102").Value = 3470.111
$ws.Range("M102").Value = -1848.111
$ws.Range("H126").Value = 2830443.5
$ws.Range("I126").Value = 3089270.5
$ws.Range("J126").Value = 501000
$ws.Range("K126").Value = 9267811.5
$ws.Range("L126").Value = 1503000
$ws.Range("M126").Value = -9265341.5
$ws.Range("N126").Value = -1507940
$ws.Range("H139").Value = 45118
$ws.Range("J139").Value = 45118
$ws.Range("L139").Value = 45118
$ws.Range("N139").Value = -55398

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 3500
$ws.Range("J11").Value = 3500
$ws.Range("L11").Value = 3500
$ws.Range("N11").Value = -3780
$ws.Range("H46").Value = 2277.9333
$ws.Range("J46").Value = 2963.2222
$ws.Range("L46").Value = 2963.2222
$ws.Range("N46").Value = -3339.2222
$ws.Range("H94").Value = 52219.332
$ws.Range("J94").Value = 52219.332
$ws.Range("L94").Value = 52219.332
$ws.Range("N94").Value = -53571.332
$ws.Range("H132").Value = 1703.7675
$ws.Range("I132").Value = 1414.238
$ws.Range("K132").Value = 4242.714
$ws.Range("M132").Value = -1712.714
$ws.Range("H138").Value = 88887.5
$ws.Range("J138").Value = 88887.5
$ws.Range("L138").Value = 88887.5
$ws.Range("N138").Value = -99167.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 39084.832
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 39084.832
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 39084.832
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -39308.832
